# Update the workbook to reflect data refresh through 2022-09-05
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (tab + internal reference) to reflect the new "through" date
$ws.Name = "Through 2022-09-05"

# Update the header label in I1 (shared string "2022 (through 09-03)" -> "2022 (through 09-05)")
$ws.Range("I1").Value = "2022 (through 09-05)"

# Update September (row 9) and October (row 10) 2022 YoY figures
$ws.Range("I9").Value = 167
$ws.Range("I10").Value = 28

# Update the Total row (row 14) to reflect the new sum
$ws.Range("I14").Value = 1166
